# Update "Pais" worksheet: refresh country COVID stats and re-sort ties
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 19 de Septiembre de 2020 a las 18:18"

# Row 4
$ws.Range("B4").Value = 6940979
$ws.Range("C4").Value = 15038
$ws.Range("D4").Value = 4193772
$ws.Range("E4").Value = 2543771
$ws.Range("G4").Value = 265
$ws.Range("H4").Value = 203436

# Row 5
$ws.Range("B5").Value = 5351723
$ws.Range("C5").Value = 46248
$ws.Range("D5").Value = 4249648
$ws.Range("E5").Value = 1016032
$ws.Range("G5").Value = 418
$ws.Range("H5").Value = 86043

# Row 6
$ws.Range("B6").Value = 4503002
$ws.Range("C6").Value = 5568
$ws.Range("E6").Value = 577828
$ws.Range("G6").Value = 178
$ws.Range("H6").Value = 136035

# Row 14
$ws.Range("B14").Value = 444674
$ws.Range("C14").Value = 1847
$ws.Range("D14").Value = 418101
$ws.Range("E14").Value = 14319
$ws.Range("G14").Value = 55
$ws.Range("H14").Value = 12254

# Row 17
$ws.Range("B17").Value = 390358
$ws.Range("C17").Value = 4422
$ws.Range("G17").Value = 27
$ws.Range("H17").Value = 41759

# Row 22
$ws.Range("B22").Value = 301348
$ws.Range("C22").Value = 1538
$ws.Range("D22").Value = 266117
$ws.Range("E22").Value = 27786
$ws.Range("G22").Value = 68
$ws.Range("H22").Value = 7445

# Row 23
$ws.Range("B23").Value = 296569
$ws.Range("C23").Value = 1638
$ws.Range("D23").Value = 217716
$ws.Range("E23").Value = 43161
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 35692

# Row 29
$ws.Range("B29").Value = 142745
$ws.Range("C29").Value = 834
$ws.Range("D29").Value = 124172
$ws.Range("E29").Value = 9362
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 9211

# Row 49
$ws.Range("B49").Value = 75461
$ws.Range("C49").Value = 231
$ws.Range("D49").Value = 73212
$ws.Range("E49").Value = 1473
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 776

# Row 57
$ws.Range("D57").Value = 57142
$ws.Range("E57").Value = 389

# Row 65
$ws.Range("B65").Value = 45877
$ws.Range("C65").Value = 117
$ws.Range("D65").Value = 45081
$ws.Range("E65").Value = 499
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 297

# Row 86
$ws.Range("B86").Value = 16557
$ws.Range("C86").Value = 140
$ws.Range("D86").Value = 13792
$ws.Range("E86").Value = 2076
$ws.Range("G86").Value = 6
$ws.Range("H86").Value = 689

# Row 88
$ws.Range("B88").Value = 14978
$ws.Range("C88").Value = 240
$ws.Range("E88").Value = 4658
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 331

# Row 93
$ws.Range("B93").Value = 12820
$ws.Range("C93").Value = 51
$ws.Range("E93").Value = 2182

# Row 112
$ws.Range("A112").Value = "Uganda"
$ws.Range("B112").Value = 6017
$ws.Range("C112").Value = 423
$ws.Range("D112").Value = 2581
$ws.Range("E112").Value = 3373
$ws.Range("G112").Value = 2
$ws.Range("H112").Value = 63

# Row 113
$ws.Range("A113").Value = "Malaui"
$ws.Range("B113").Value = 5716
$ws.Range("D113").Value = 4026
$ws.Range("E113").Value = 1511
$ws.Range("H113").Value = 179

# Row 124
$ws.Range("A124").Value = "Jamaica"
$ws.Range("B124").Value = 4758
$ws.Range("C124").Value = 187
$ws.Range("D124").Value = 1327
$ws.Range("E124").Value = 3371
$ws.Range("G124").Value = 5
$ws.Range("H124").Value = 60

# Row 125
$ws.Range("A125").Value = "Surinam"
$ws.Range("B125").Value = 4691
$ws.Range("D125").Value = 4280
$ws.Range("E125").Value = 315
$ws.Range("H125").Value = 96

# Row 126
$ws.Range("A126").Value = "Ruanda"
$ws.Range("B126").Value = 4671
$ws.Range("D126").Value = 2845
$ws.Range("E126").Value = 1801
$ws.Range("H126").Value = 25

# Row 127
$ws.Range("B127").Value = 4540
$ws.Range("C127").Value = 196
$ws.Range("D127").Value = 2672
$ws.Range("E127").Value = 1838
$ws.Range("G127").Value = 1
$ws.Range("H127").Value = 30

# Row 151
$ws.Range("B151").Value = 2159
$ws.Range("C151").Value = 6
$ws.Range("D151").Value = 1650
$ws.Range("E151").Value = 437

# Row 159
$ws.Range("B159").Value = 1580
$ws.Range("C159").Value = 15
$ws.Range("E159").Value = 276

# Row 163
$ws.Range("B163").Value = 1335
$ws.Range("C163").Value = 1
$ws.Range("D163").Value = 1216
$ws.Range("E163").Value = 37

# Row 204
$ws.Range("A204").Value = "Santa Lucia"

# Row 205
$ws.Range("A205").Value = "Timor Oriental"

# Row 214
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Row 215
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0
